# Auto-generated edit script applying the Golem_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1333
$ws.Range("J2").Value = 1500
$ws.Range("L2").Value = 1500
$ws.Range("N2").Value = -1726
$ws.Range("H6").Value = 81.92308
$ws.Range("I6").Value = 105
$ws.Range("K6").Value = 315
$ws.Range("M6").Value = -203
$ws.Range("H38").Value = 27.222221
$ws.Range("I38").Value = 27.222221
$ws.Range("K38").Value = 81.666663
$ws.Range("M38").Value = 290.333337
$ws.Range("H39").Value = 6839.9375
$ws.Range("I39").Value = 814.8333
$ws.Range("K39").Value = 2444.4999
$ws.Range("M39").Value = -2148.4999
$ws.Range("H41").Value = 687.25
$ws.Range("J41").Value = 979.5
$ws.Range("L41").Value = 979.5
$ws.Range("N41").Value = -1859.5
$ws.Range("H93").Value = 750089500
$ws.Range("J93").Value = 36500
$ws.Range("L93").Value = 36500
$ws.Range("N93").Value = -41492
$ws.Range("H99").Value = 83333544
$ws.Range("I99").Value = 83333544
$ws.Range("K99").Value = 250000632
$ws.Range("M99").Value = -249999134
$ws.Range("H106").Value = 250000430
$ws.Range("I106").Value = 250000430
$ws.Range("K106").Value = 250000430
$ws.Range("M106").Value = -249999799
$ws.Range("H107").Value = 91670
$ws.Range("I107").Value = 101522.22
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 101522.22
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -99602.22
$ws.Range("N107").Value = -6840
$ws.Range("H123").Value = 130000
$ws.Range("J123").Value = 130000
$ws.Range("L123").Value = 130000
$ws.Range("N123").Value = -139800
$ws.Range("H138").Value = 7554
$ws.Range("J138").Value = 7205.3335
$ws.Range("L138").Value = 21616.0005
$ws.Range("N138").Value = -31896.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13989.8
$ws.Range("I32").Value = 12487.25
$ws.Range("K32").Value = 12487.25
$ws.Range("M32").Value = -12200.25
$ws.Range("H101").Value = 31900
$ws.Range("J101").Value = 31900
$ws.Range("L101").Value = 31900
$ws.Range("N101").Value = -38390
$ws.Range("H110").Value = 125000400
$ws.Range("I110").Value = 799
$ws.Range("K110").Value = 799
$ws.Range("M110").Value = 1246

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 101112.25
$ws.Range("I107").Value = 200725
$ws.Range("J107").Value = 1499.5
$ws.Range("K107").Value = 200725
$ws.Range("L107").Value = 1499.5
$ws.Range("M107").Value = -198805
$ws.Range("N107").Value = -5339.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3011.4443
$ws.Range("I105").Value = 2220.6
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 2220.6
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -473.5999999999999
$ws.Range("N105").Value = -7494
$ws.Range("H132").Value = 1374.3334
$ws.Range("J132").Value = 2100
$ws.Range("L132").Value = 6300
$ws.Range("N132").Value = -11360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 300
$ws.Range("J9").Value = 316.66666
$ws.Range("L9").Value = 949.9999799999999
$ws.Range("N9").Value = -1397.99998
$ws.Range("H10").Value = 43
$ws.Range("I10").Value = 21.238094
$ws.Range("J10").Value = 500
$ws.Range("K10").Value = 63.714282
$ws.Range("L10").Value = 1500
$ws.Range("M10").Value = 75.285718
$ws.Range("N10").Value = -1778
$ws.Range("H11").Value = 30
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").Value = ""
$ws.Range("H13").Value = 750.2
$ws.Range("I13").Value = 816.3333
$ws.Range("J13").Value = 651
$ws.Range("K13").Value = 2448.9999
$ws.Range("L13").Value = 1953
$ws.Range("M13").Value = -2280.9999
$ws.Range("N13").Value = -2289
$ws.Range("H16").Value = 10
$ws.Range("J16").Value = 10
$ws.Range("L16").Value = 30
$ws.Range("N16").Value = -376
$ws.Range("H17").Value = 1000
$ws.Range("I17").Value = 1000
$ws.Range("K17").Value = 3000
$ws.Range("M17").Value = -2831
$ws.Range("H46").Value = 2246.647
$ws.Range("I46").Value = 399.33334
$ws.Range("J46").Value = 2642.5
$ws.Range("K46").Value = 1198.00002
$ws.Range("L46").Value = 7927.5
$ws.Range("M46").Value = -1107.00002
$ws.Range("N46").Value = -8109.5
$ws.Range("H80").Value = 2180.8
$ws.Range("J80").Value = 2202
$ws.Range("L80").Value = 6606
$ws.Range("N80").Value = -8478
$ws.Range("H83").Value = 2180.8
$ws.Range("J83").Value = 2202
$ws.Range("L83").Value = 19818
$ws.Range("N83").Value = -29178
$ws.Range("H115").Value = 2749.5
$ws.Range("I115").Value = 1500
$ws.Range("J115").Value = 3166
$ws.Range("K115").Value = 4500
$ws.Range("L115").Value = 9498
$ws.Range("M115").Value = -3325
$ws.Range("N115").Value = -11848
$ws.Range("H117").Value = 203.5
$ws.Range("I117").Value = 203.5
$ws.Range("K117").Value = 610.5
$ws.Range("M117").Value = 2831.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 997.5
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = ""
$ws.Range("H68").Value = 5000
$ws.Range("I68").Value = 5000
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 5000
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -4251
$ws.Range("N68").Value = ""
$ws.Range("H71").Value = 5000
$ws.Range("I71").Value = 5000
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 25000
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -21256
$ws.Range("N71").Value = ""
$ws.Range("H93").Value = 17548970
$ws.Range("I93").Value = 22226850
$ws.Range("K93").Value = 22226850
$ws.Range("M93").Value = -22225602

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2124.75
$ws.Range("J132").Value = 2500
$ws.Range("L132").Value = 7500
$ws.Range("N132").Value = -12560
